$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Test Cases")
$src = $ws1.Range("D31")
$src.Copy()
foreach ($r in @(29,30,32,33,34,35,36,37,38,39)) {
    $ws1.Range("D$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
$wb.Save()
